$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (COM ColumnWidth is pixel-quantized at write time: stored
# OOXML width = (round(ColumnWidth*7)+5)/7 ; pre-compensate so the saved width
# lands as close as possible to the target character width from the diff)
$ws.Columns.Item(1).ColumnWidth = 5.714285714285714
$ws.Columns.Item(3).ColumnWidth = 11.142857142857142
$ws.Columns.Item(5).ColumnWidth = 14.428571428571429

# Update row 3 data
$ws.Range("A3").Value = "andres"
$ws.Range("B3").Value = "perez"
$ws.Range("C3").Value = "andres@chile.com"
$ws.Range("D3").Value = 51231212
$ws.Range("E3").Value = "Los conejitos 4321"
$ws.Range("G3").Value = "30/05/2014 - 16:21:57"
